# 06/11/2025 Fixed Showing up of Draft Button for IT PIC (Accepted Tix)
# Adds six new trailing columns (P:U) to the ClosedTicket header row (row 3):
#   ASSIGNED IT PIC / ASSIGNED DATE TIME / RESOLVED DATE TIME / SLA HOURS / ACTUAL HOURS / HIT OR MISS
# and widens the new columns to match the author's saved column widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 3, columns P:U), values taken from the commit's new shared strings ---
$ws.Range("P3").Value = "ASSIGNED IT PIC "
$ws.Range("Q3").Value = "ASSIGNED DATE TIME"
$ws.Range("R3").Value = "RESOLVED DATE TIME "
$ws.Range("S3").Value = "SLA HOURS"
$ws.Range("T3").Value = "ACTUAL HOURS"
$ws.Range("U3").Value = "HIT OR MISS"

# Match the existing header formatting (fill/style) used by A3:O3 - copy format from O3 onto the new cells
$ws.Range("O3").Copy()
$ws.Range("P3:U3").PasteSpecial(-4122)  # xlPasteFormats

# --- Column widths for the new columns, matching the widths saved in the workbook ---
$ws.Columns.Item(16).ColumnWidth = 34.833333333333336  # P -> 35.7109375
$ws.Columns.Item(17).ColumnWidth = 37.666666666666664  # Q -> 38.5703125
$ws.Columns.Item(18).ColumnWidth = 38.666666666666664  # R -> 39.42578125
$ws.Columns.Item(19).ColumnWidth = 36.333333333333336  # S -> 37.140625
$ws.Columns.Item(20).ColumnWidth = 34.833333333333336  # T -> 35.7109375
$ws.Columns.Item(21).ColumnWidth = 18.666666666666668  # U -> 19.5703125

# --- View state: scrolled right so the new columns are visible, selection on the last new header cell ---
$excel.ActiveWindow.Zoom = 70
$excel.ActiveWindow.ScrollColumn = 15
$null = $ws.Range("U3").Select()
